# Fruta / hortaliza, semanal
#
# A new weekly price record is inserted as row 80 (pushing all the
# subsequent records down by one row, so the former row 131 becomes
# row 132). The new row re-uses the same market/category/quality
# metadata as the record that used to sit in row 80 (now shifted to
# row 81), only the date (column D) and volume (column J) differ.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 80; everything from row 80 down
# (through the former last row, 131) shifts down to make room.
$ws.Rows.Item(80).Insert()

# The record that used to be row 80 is now row 81. Duplicate it into
# the newly-created row 80 so every column besides the date/volume
# keeps the same values.
$ws.Range("A81:R81").Copy()
$ws.Range("A80").PasteSpecial()

# Apply the new record's own date and volume.
$ws.Range("D80").Value2 = 44873
$ws.Range("J80").Value2 = 1000
